{"js": "// Applies the replacements described by the diff: the date line plus the\n// 25 \"NNN\u00f7N=\" table-cell expressions are swapped for their updated values.\n// Each mapping is an exact, unique old->new text pair (verified against the\n// source document), so a body-wide search+replace for each pair reproduces\n// the diff exactly without touching any other run formatting.\nconst replacements = [\n  [\"2024-03-18 Monday\", \"2024-03-19 Tuesday\"],\n  [\"718\u00f72=\", \"945\u00f78=\"],\n  [\"350\u00f77=\", \"844\u00f74=\"],\n  [\"394\u00f78=\", \"705\u00f79=\"],\n  [\"219\u00f74=\", \"954\u00f75=\"],\n  [\"386\u00f76=\", \"278\u00f78=\"],\n  [\"998\u00f78=\", \"252\u00f78=\"],\n  [\"616\u00f78=\", \"366\u00f74=\"],\n  [\"829\u00f73=\", \"322\u00f75=\"],\n  [\"659\u00f75=\", \"934\u00f78=\"],\n  [\"661\u00f77=\", \"800\u00f79=\"],\n  [\"727\u00f77=\", \"556\u00f72=\"],\n  [\"569\u00f74=\", \"164\u00f74=\"],\n  [\"110\u00f75=\", \"594\u00f79=\"],\n  [\"539\u00f72=\", \"748\u00f74=\"],\n  [\"102\u00f75=\", \"317\u00f75=\"],\n  [\"740\u00f77=\", \"666\u00f75=\"],\n  [\"943\u00f76=\", \"413\u00f75=\"],\n  [\"699\u00f73=\", \"961\u00f77=\"],\n  [\"573\u00f73=\", \"994\u00f72=\"],\n  [\"523\u00f76=\", \"253\u00f73=\"],\n  [\"117\u00f72=\", \"463\u00f79=\"],\n  [\"392\u00f77=\", \"470\u00f72=\"],\n  [\"941\u00f74=\", \"426\u00f73=\"],\n  [\"941\u00f73=\", \"373\u00f73=\"],\n  [\"571\u00f74=\", \"847\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the replacements described by the diff: the date line plus the\n# 25 \"NNN\u00f7N=\" table-cell expressions are swapped for their updated values.\n# Each mapping is an exact, unique old->new text pair, so a document-wide\n# Find/Replace (wdReplaceAll) for each pair reproduces the diff exactly\n# without disturbing any other run formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-03-18 Monday\"; New = \"2024-03-19 Tuesday\" }\n    @{ Old = \"718\u00f72=\"; New = \"945\u00f78=\" }\n    @{ Old = \"350\u00f77=\"; New = \"844\u00f74=\" }\n    @{ Old = \"394\u00f78=\"; New = \"705\u00f79=\" }\n    @{ Old = \"219\u00f74=\"; New = \"954\u00f75=\" }\n    @{ Old = \"386\u00f76=\"; New = \"278\u00f78=\" }\n    @{ Old = \"998\u00f78=\"; New = \"252\u00f78=\" }\n    @{ Old = \"616\u00f78=\"; New = \"366\u00f74=\" }\n    @{ Old = \"829\u00f73=\"; New = \"322\u00f75=\" }\n    @{ Old = \"659\u00f75=\"; New = \"934\u00f78=\" }\n    @{ Old = \"661\u00f77=\"; New = \"800\u00f79=\" }\n    @{ Old = \"727\u00f77=\"; New = \"556\u00f72=\" }\n    @{ Old = \"569\u00f74=\"; New = \"164\u00f74=\" }\n    @{ Old = \"110\u00f75=\"; New = \"594\u00f79=\" }\n    @{ Old = \"539\u00f72=\"; New = \"748\u00f74=\" }\n    @{ Old = \"102\u00f75=\"; New = \"317\u00f75=\" }\n    @{ Old = \"740\u00f77=\"; New = \"666\u00f75=\" }\n    @{ Old = \"943\u00f76=\"; New = \"413\u00f75=\" }\n    @{ Old = \"699\u00f73=\"; New = \"961\u00f77=\" }\n    @{ Old = \"573\u00f73=\"; New = \"994\u00f72=\" }\n    @{ Old = \"523\u00f76=\"; New = \"253\u00f73=\" }\n    @{ Old = \"117\u00f72=\"; New = \"463\u00f79=\" }\n    @{ Old = \"392\u00f77=\"; New = \"470\u00f72=\" }\n    @{ Old = \"941\u00f74=\"; New = \"426\u00f73=\" }\n    @{ Old = \"941\u00f73=\"; New = \"373\u00f73=\" }\n    @{ Old = \"571\u00f74=\"; New = \"847\u00f76=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
